# Generate Report for Handoff
# Updates the "f72e2c20..." handoff entry to the new "6ea0face..." handoff
# (new .md commit + new translation package), and adds a brand-new handoff
# row for "ffff5a5d90a6-3c5e-450a-bd40-6e9a85502d56.md" on every sheet.

$wb = $excel.ActiveWorkbook

$oldUuid = "f72e2c20-d5e3-4a4e-8ae9-080c60ea4178"
$newUuid = "6ea0face-9cb0-45d6-8b7d-dfac503be676"
$oldHash = "492faa3d70de96bc2f8e5b0ca3c787164bf033b3"
$newHash = "d072f7e8a3d914fadbfd5d199cad32120bbf389b"
$newFileUuid = "ffff5a5d90a6-3c5e-450a-bd40-6e9a85502d56"

$newMdName = "$newUuid.md"
$newFileMdName = "$newFileUuid.md"
$newZhXlfName = "$newUuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newUuid.$newHash.de-de.xlf"

$newHandoffDate = "2016-03-18 17:13:41"
$newHandoffDatetime = "2016-03-18 17:13:33"

$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/f8d0b906c4464ea3cd9c5056d102ff096481255a/e2e"
$zhUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1fd3a1019b694eb09c9d7e424c5346ee0cddb93/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht"
$deUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/87e9e5cd5eb0a743f606ef192604dc199cbdb1ed/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht"

$newMdUrl = "$mdUrlBase/$newMdName"
$newFileMdUrl = "$mdUrlBase/$newFileMdName"
$newZhUrl = "$zhUrlBase/$newZhXlfName"
$newDeUrl = "$deUrlBase/$newDeXlfName"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Update existing row 2 in place: new md name + new handoff date
$ov.Range("A2").Value = $newMdName
$ov.Range("D2").Value = $newHandoffDate

# Append new row 3 for the new handoff file
$ov.Range("A3").Value = $newFileMdName
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = $newHandoffDate
$ov.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Range("A3"), $newFileMdUrl, "", "", $newFileMdName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" and "de-de": both share the same 12-column layout
#   A Source File Name | B File Extension | C Status | D Latest Handoff File
#   E Latest Handoff Datetime | H Latest Handback DateTime | I Reference Tokens
#   J Handoff Reason
# ---------------------------------------------------------------------
function Update-LangSheet($sheetName, $targetXlfName, $targetXlfUrl) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Update existing row 2 in place
    $ws.Range("A2").Value = $newMdName
    $ws.Range("D2").Value = $targetXlfName
    $ws.Range("E2").Value = $newHandoffDatetime

    # Append new row 3, mirroring row 2 but with the new handoff file name
    $ws.Range("A3").Value = $newFileMdName
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $targetXlfName
    $ws.Range("E3").Value = $newHandoffDatetime
    $ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("I3").Value = ""
    $ws.Range("J3").Value = "Include"

    $ws.Hyperlinks.Add($ws.Range("A3"), $newFileMdUrl, "", "", $newFileMdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), $newFileMdUrl, "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $targetXlfUrl, "", "", $targetXlfName) | Out-Null
}

Update-LangSheet "zh-cn" $newZhXlfName $newZhUrl
Update-LangSheet "de-de" $newDeXlfName $newDeUrl

Write-Output "edit complete"
